$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Private" column (D) as TRUE for all data rows (2-6)
$ws.Range("D2:D6").Value = $true

# Column D now also accepts the TRUE/FALSE list (same validation rule as column F)
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("D2:D6").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Reflect the new active selection left after the edit
[void]$ws.Range("D2:D6").Select()
